# Delete the "adductName" column (column C) from the "Corrected" sheet.
# This shifts columns D:G left to C:F, so the sheet ends up with data in
# A:F instead of A:G (matches the xl/worksheets/sheet2.xml diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")

$ws.Columns.Item(3).EntireColumn.Delete()
